$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.027783513069153
$ws.Range("B1").Value = 1.638083219528198
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.729069948196411
$ws.Range("E1").Value = 1.351328611373901
